$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Customer Class section (rows 3-6): fill in "Total Points" (column E) ---
# to match full marks awarded for each sub-item (row 7 auto-sums via formula)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# --- Product Class section (rows 10-14): fill in "Total Points" (column E) ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2

# toString() method (row 14) loses 1 point for incorrect format, with comment
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = "(-1) for incorrect format for toString method"

# Move the active selection to F12 (no more scrolled-down topLeftCell)
[void]$ws.Range("F12").Select()
